$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two earned-value cells for rows 8 and 9 (Core Features section)
$ws.Range("C8").Value = 0.1
$ws.Range("C9").Value = 0

# Update the view state: scroll so row 7 is the top-left visible row,
# and select C9 as the active cell
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C9").Select()
